# Appends 24 new gesture-collection rows (Josh136-Josh159) to Sheet1,
# reconciling the "Laying" and "Running" sample sets as described in the
# commit message ("reconciled duplicates and file structure").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    ,@('Josh136', 'Laying down, hands at rest on chest, when a subtle tonic and tremor episode occurs halfway through', 'Laying', 'Hands at side/rest', 'Tremor', 'Laying (hands at rest)', 'Laying (hands tremor)')
    ,@('Josh137', 'Laying down, hands at rest on chest, when a strong tonic and tremor episode occurs halfway through', 'Laying', 'Hands at side/rest', 'Tremor', 'Laying (hands at rest)', 'Laying (hands tremor)')
    ,@('Josh138', 'Laying down, hands at rest on chest, when a extreme tonic and tremor episode occurs halfway through', 'Laying', 'Hands at side/rest', 'Tremor', 'Laying (hands at rest)', 'Laying (hands tremor)')
    ,@('Josh139', 'Laying down, using mobile phone (swiping), when a subtle tonic and tremor episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tremor', 'Laying (hands active)', 'Laying (hands tremor)')
    ,@('Josh140', 'Laying down, using mobile phone (swiping), when a subtle tonic and tremor episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tremor', 'Laying (hands active)', 'Laying (hands tremor)')
    ,@('Josh141', 'Laying down, using mobile phone (swiping), when a subtle tonic and tremor episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tremor', 'Laying (hands active)', 'Laying (hands tremor)')
    ,@('Josh142', 'Laying down, using remote and watching TV, when a subtle tonic and tremor episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tremor', 'Laying (hands active)', 'Laying (hands tremor)')
    ,@('Josh143', 'Laying down, using remote and watching TV, when a subtle tonic and tremor episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tremor', 'Laying (hands active)', 'Laying (hands tremor)')
    ,@('Josh144', 'Laying down, using remote and watching TV, when a subtle tonic and tremor episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tremor', 'Laying (hands active)', 'Laying (hands tremor)')
    ,@('Josh145', 'Laying down, hands at rest on chest, when a subtle tonic episode occurs halfway through', 'Laying', 'Hands at side/rest', 'Tonic', 'Laying (hands at rest)', 'Laying (hands tonic)')
    ,@('Josh146', 'Laying down, hands at rest on chest, when a strong tonic episode occurs halfway through', 'Laying', 'Hands at side/rest', 'Tonic', 'Laying (hands at rest)', 'Laying (hands tonic)')
    ,@('Josh147', 'Laying down, hands at rest on chest, when a extreme tonic episode occurs halfway through', 'Laying', 'Hands at side/rest', 'Tonic', 'Laying (hands at rest)', 'Laying (hands tonic)')
    ,@('Josh148', 'Laying down, using mobile phone (swiping), when a subtle tonic episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tonic', 'Laying (hands active)', 'Laying (hands tonic)')
    ,@('Josh149', 'Laying down, using mobile phone (swiping), when a strong tonic episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tonic', 'Laying (hands active)', 'Laying (hands tonic)')
    ,@('Josh150', 'Laying down, using mobile phone (swiping), when a extreme tonic episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tonic', 'Laying (hands active)', 'Laying (hands tonic)')
    ,@('Josh151', 'Laying down, using remote and watching TV, when a subtle tonic episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tonic', 'Laying (hands active)', 'Laying (hands tonic)')
    ,@('Josh152', 'Laying down, using remote and watching TV, when a strong tonic episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tonic', 'Laying (hands active)', 'Laying (hands tonic)')
    ,@('Josh153', 'Laying down, using remote and watching TV, when a extreme tonic episode occurs halfway through', 'Laying', 'Hands active (moving)', 'Tonic', 'Laying (hands active)', 'Laying (hands tonic)')
    ,@('Josh154', 'Running down the street, hands at side swinging, and then stopping when a subtle tonic episode occurs halfway through', 'Running', 'Hands at rest', 'Tonic', 'Running (hands at rest)', 'Running (hands tremor)')
    ,@('Josh155', 'Running down the street, hands at side swinging, and then stopping when a strong tonic episode occurs halfway through', 'Running', 'Hands at rest', 'Tonic', 'Running (hands at rest)', 'Running (hands tremor)')
    ,@('Josh156', 'Running down the street, hands at side swinging, and then stopping when a extreme tonic episode occurs halfway through', 'Running', 'Hands at rest', 'Tonic', 'Running (hands at rest)', 'Running (hands tremor)')
    ,@('Josh157', 'Running down the street, hands held down at sides (minimal movement), and then stopping when a subtle tonic episode occurs halfway through', 'Running', 'Hands at side/rest', 'Tonic', 'Running (hands at rest)', 'Running (hands tremor)')
    ,@('Josh158', 'Running down the street, hands held down at sides (minimal movement), and then stopping when a strong tonic episode occurs halfway through', 'Running', 'Hands at side/rest', 'Tonic', 'Running (hands at rest)', 'Running (hands tremor)')
    ,@('Josh159', 'Running down the street, hands held down at sides (minimal movement), and then stopping when a extreme tonic episode occurs halfway through', 'Running', 'Hands at side/rest', 'Tonic', 'Running (hands at rest)', 'Running (hands tremor)')
)

$startRow = 137
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

$ws.Range("B149").Select()
